$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new detail rows (2509 period) right after the existing
#     2508 rows (16-17), pushing the signature block (old rows 22-23) down
#     to rows 24-25, matching the target layout (dimension B2:J25). ---
$ws.Rows("18:19").Insert()

# Copy the formatting (borders/fills/fonts) + values from the 2508 rows
# (16-17) down into the freshly inserted rows so the new entries keep the
# same bordered-table look as the rest of "BD" table (borderId group
# matching rows 16/17, with row 19 picking up the closing/bottom border
# that used to belong to row 17).
$ws.Range("B16:J17").Copy($ws.Range("B18:J19"))
$excel.CutCopyMode = 0

# --- New row 18: YOIDER RODRIGUEZ BLANCO, periodo 2509 ---
$ws.Range("B18").Value2 = "CC"
$ws.Range("C18").Value2 = "1002320130"
$ws.Range("D18").Value2 = "YOIDER RODRIGUEZ BLANCO"
$ws.Range("E18").Value2 = "2509"
$ws.Range("F18").Value2 = 56940
$ws.Range("G18").Value2 = 1423500
$ws.Range("H18").Value2 = ""
$ws.Range("I18").Value2 = ""
$ws.Range("J18").Value2 = ""

# --- New row 19: KEINER ENRIQUE TORRES GOMEZ, periodo 2509 ---
$ws.Range("B19").Value2 = "CC"
$ws.Range("C19").Value2 = "1007972996"
$ws.Range("D19").Value2 = "KEINER ENRIQUE TORRES GOMEZ"
$ws.Range("E19").Value2 = "2509"
$ws.Range("F19").Value2 = 56940
$ws.Range("G19").Value2 = 1423500
$ws.Range("H19").Value2 = ""
$ws.Range("I19").Value2 = ""
$ws.Range("J19").Value2 = ""

# --- Totals: one more "periodo" and double the accrued "VALOR MORA" now
#     that both workers show up for 2508 and 2509. ---
$ws.Range("E11").Value2 = 227760
$ws.Range("F13").Value2 = 2
